# Insert a new data row before the current row 60, shifting existing rows
# 60-92 down to 61-93 (new weekly observation added to the series).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(60).Insert()

$ws.Range("A60").Value = 6
$ws.Range("B60").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C60").Value = "Metropolitana"
$ws.Range("D60").Value = 45119
$ws.Range("E60").Value = 13
$ws.Range("F60").Value = 100112035
$ws.Range("G60").Value = "Bruselas (repollito)"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 157
$ws.Range("K60").Value = 17000
$ws.Range("L60").Value = 18000
$ws.Range("M60").Value = 17236
$ws.Range("N60").Value = "$/malla 15 kilos"
$ws.Range("O60").Value = "Provincia de Quillota"
$ws.Range("P60").Value = 1149
$ws.Range("Q60").Value = 15
$ws.Range("R60").Value = "Hortaliza"
